$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: D-column price values. Many look like plain decimals (e.g. "554.92")
# which Excel's COM layer would otherwise auto-coerce to a Number; prefixing
# with an apostrophe keeps them stored as Text, matching the source data
# (multi-dot, thousands-style values like "63.834.85" are never ambiguous and
# don't need the prefix, but it is harmless to add it everywhere).
function Set-D($row, $val) {
    $ws.Cells.Item($row, 4).Value = "'" + $val
}
function Set-E($row, $val) {
    $ws.Cells.Item($row, 5).Value = "  $val  "
}
function Set-B($row, $val) {
    $ws.Cells.Item($row, 2).Value = $val
}
function Set-C($row, $val) {
    $ws.Cells.Item($row, 3).Value = $val
}

# Row 2 - Bitcoin
Set-D 2 "63.834.85"
Set-E 2 "-2.56%"

# Row 3 - Ethereum
Set-D 3 "3.336.02"
Set-E 3 "-3.87%"

# Row 4 - TetherUSD
Set-E 4 "+0.04%"

# Row 5 - BNB
Set-D 5 "554.92"

# Row 6 - Solana
Set-D 6 "175.05"
Set-E 6 "-0.28%"

# Row 7 - XRP
Set-E 7 "-2.66%"

# Row 8 - LidoStakedEther
Set-D 8 "3.329.36"
Set-E 8 "-4.03%"

# Row 9 - USDC
Set-E 9 "+0.07%"

# Row 10 - Cardano
Set-D 10 "0.620"
Set-E 10 "-0.53%"

# Row 11 - Dogecoin
Set-E 11 "+4.38%"

# Row 12 - Avalanche
Set-D 12 "54.07"
Set-E 12 "+1.53%"

# Row 13 - ShibaInu
Set-E 13 "+0.77%"

# Row 14 - Polkadot
Set-D 14 "9.00"
Set-E 14 "-0.65%"

# Row 15 - WrappedliquidstakedEther2.0
Set-D 15 "3.870.80"
Set-E 15 "-3.70%"

# Row 16 - Chainlink
Set-D 16 "18.29"
Set-E 16 "+1.13%"

# Row 17 - TRON
Set-E 17 "-2.52%"

# Row 18 - WrappedEther
Set-D 18 "3.337.97"
Set-E 18 "-3.84%"

# Row 19 - Uniswap
Set-D 19 "11.82"
Set-E 19 "-1.04%"

# Row 20 - WrappedBTC
Set-D 20 "63.777.50"
Set-E 20 "-2.56%"

# Row 21 - Polygon
Set-D 21 "0.974"
Set-E 21 "-1.22%"

# Row 22 - BitcoinCash
Set-D 22 "429.10"
Set-E 22 "+4.61%"

# Row 23 - Toncoin
Set-E 23 "+9.94%"

# Row 24 - PancakeSwap
Set-D 24 "4.11"
Set-E 24 "+0.24%"

# Row 25 - Litecoin
Set-D 25 "83.97"
Set-E 25 "-0.76%"

# Row 26 - InternetComputer(DFINITY)
Set-D 26 "13.05"
Set-E 26 "+4.01%"

# Row 27 - RenderToken
Set-D 27 "10.61"
Set-E 27 "-1.70%"

# Row 28 - ImmutableX
Set-D 28 "2.81"
Set-E 28 "+0.49%"

# Row 29 - Filecoin
Set-E 29 "-1.79%"

# Row 30 - EthereumClassic
Set-D 30 "29.52"
Set-E 30 "-1.22%"

# Row 31 - NEARProtocol
Set-D 31 "6.52"
Set-E 31 "+4.08%"

# Row 32 - Bittensor
Set-D 32 "589.96"
Set-E 32 "-3.85%"

# Row 33 - Cosmos
Set-D 33 "11.40"
Set-E 33 "-0.86%"

# Row 34 - Hedera
Set-E 34 "-1.16%"

# Row 35 - OKB
Set-D 35 "58.54"
Set-E 35 "-0.46%"

# Row 36 - Dai
Set-E 36 "+0.10%"

# Row 37 - Kaspa
Set-E 37 "-4.51%"

# Row 38 - Stacks
Set-E 38 "+4.28%"

# Row 39 - InjectiveProtocol
Set-E 39 "-3.13%"

# Row 40 - PEPE (note: contains U+2083 SUBSCRIPT THREE, written literally below)
Set-D 40 "0.0₃0748"
Set-E 40 "-4.14%"

# Row 41 - TheGraph
Set-E 41 "-2.63%"

# Row 42 - Maker
Set-D 42 "3.108.29"
Set-E 42 "-5.75%"

# Row 43 - FirstDigitalUSD
Set-E 43 "-0.09%"

# Row 44 - ThetaToken
Set-D 44 "2.82"
Set-E 44 "-0.75%"

# Row 45 - was VeChain, now ApeXProtocol (rows 45/46 swapped places)
Set-B 45 "ApeXProtocol"
Set-C 45 "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-D 45 "3.18"
Set-E 45 "-1.87%"

# Row 46 - was ApeXProtocol, now VeChain (rows 45/46 swapped places)
Set-B 46 "VeChain"
Set-C 46 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-D 46 "0.0405"
Set-E 46 "-1.03%"

# Row 47 - Fetch.AI
Set-D 47 "2.43"
Set-E 47 "-1.94%"

# Row 48 - Stellar
Set-E 48 "-1.45%"

# Row 49 - WEMIXToken
Set-E 49 "-3.24%"

# Row 50 - Monero
Set-D 50 "135.12"
Set-E 50 "-1.57%"

# Row 51 - THORChain
Set-D 51 "8.17"
Set-E 51 "-1.72%"
